$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the cell values first.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin border all around, centered / top-aligned.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160

# Copy B1's format onto A2 so both share the exact same style record
# instead of building a second, near-duplicate style.
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)
